# Modified the preferential queue. The request is shifted to the left to
# get additional free space on the right.
#
# Adds three new "servers" (s4, s5, s6) to the service/request tables,
# updates the existing queue numbers, and rescales the sending interval.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "mec": collapse the per-machine queue counters down to 1.
# ---------------------------------------------------------------------
$wsMec = $wb.Worksheets.Item("mec")
[void]$wsMec.Activate()
$wsMec.Range("B2").Value = 1
$wsMec.Range("B3").Value = 1
$wsMec.Range("B4").Value = 1
[void]$wsMec.Range("B5").Select()

# ---------------------------------------------------------------------
# Sheet "service": update existing services, add s4/s5/s6 rows.
# ---------------------------------------------------------------------
$wsService = $wb.Worksheets.Item("service")
[void]$wsService.Activate()

$wsService.Range("B2").Value = 90
$wsService.Range("C2").Value = 3000
$wsService.Range("B3").Value = 22
$wsService.Range("C3").Value = 3000
$wsService.Range("B4").Value = 10
$wsService.Range("C4").Value = 3000

$wsService.Range("A5").Value = "s4"
$wsService.Range("B5").Value = 90
$wsService.Range("C5").Value = 2000

$wsService.Range("A6").Value = "s5"
$wsService.Range("B6").Value = 22
$wsService.Range("C6").Value = 2000

$wsService.Range("A7").Value = "s6"
$wsService.Range("B7").Value = 10
$wsService.Range("C7").Value = 2000

$wsService.Range("C2:C7").Font.Bold = $true
$wsService.Range("A2:C7").RowHeight = 13.8

[void]$wsService.Range("C5").Select()

# ---------------------------------------------------------------------
# Sheet "request": add s4/s5/s6 columns (E, F, G) and update values.
# ---------------------------------------------------------------------
$wsRequest = $wb.Worksheets.Item("request")
[void]$wsRequest.Activate()

$wsRequest.Range("E1").Value = "s4"
$wsRequest.Range("F1").Value = "s5"
$wsRequest.Range("G1").Value = "s6"

$wsRequest.Range("B2").Value = 700
$wsRequest.Range("C2").Value = 500
$wsRequest.Range("D2").Value = 200
$wsRequest.Range("E2").Value = 700
$wsRequest.Range("F2").Value = 500
$wsRequest.Range("G2").Value = 200

$wsRequest.Range("B3").Value = 200
$wsRequest.Range("C3").Value = 400
$wsRequest.Range("D3").Value = 700
$wsRequest.Range("E3").Value = 200
$wsRequest.Range("F3").Value = 400
$wsRequest.Range("G3").Value = 700

$wsRequest.Range("B4").Value = 400
$wsRequest.Range("C4").Value = 700
$wsRequest.Range("D4").Value = 200
$wsRequest.Range("E4").Value = 400
$wsRequest.Range("F4").Value = 700
$wsRequest.Range("G4").Value = 200

$wsRequest.Range("A1:G4").RowHeight = 13.8

[void]$wsRequest.Range("E2").Select()

# ---------------------------------------------------------------------
# Sheet "intervalForSendingRequests": scale the interval value up.
# ---------------------------------------------------------------------
$wsInterval = $wb.Worksheets.Item("intervalForSendingRequests")
[void]$wsInterval.Activate()

$wsInterval.Range("B1").Value = 120000
$wsInterval.Range("B1").Font.Bold = $true
$wsInterval.Range("A1:B1").RowHeight = 13.8

[void]$wsInterval.Range("C9").Select()
